$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range('D2')
$c.NumberFormat = '@'
$c.Value = '25.806.35'
$c.ClearFormats()
$c = $ws.Range('D3')
$c.NumberFormat = '@'
$c.Value = '1.634.53'
$c.ClearFormats()
$ws.Range('E3').Value = '  -0.21%  '
$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '215.50'
$c.ClearFormats()
$ws.Range('E5').Value = '  +0.16%  '
$ws.Range('E6').Value = '  -0.90%  '
$ws.Range('E7').Value = '  -0.25%  '
$ws.Range('E8').Value = '  +0.28%  '
$ws.Range('E9').Value = '  -0.98%  '
$c = $ws.Range('D10')
$c.NumberFormat = '@'
$c.Value = '19.67'
$c.ClearFormats()
$ws.Range('E10').Value = '  -3.28%  '
$ws.Range('E11').Value = '  +1.45%  '
$ws.Range('E12').Value = '  -0.19%  '
$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$c = $ws.Range('D13')
$c.NumberFormat = '@'
$c.Value = '1.636.22'
$c.ClearFormats()
$ws.Range('E13').Value = '  -0.02%  '
$ws.Range('B14').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C14').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$c = $ws.Range('D14')
$c.NumberFormat = '@'
$c.Value = '1.859.49'
$c.ClearFormats()
$ws.Range('E14').Value = '  -0.23%  '
$ws.Range('E15').Value = '  -0.15%  '
$c = $ws.Range('D16')
$c.NumberFormat = '@'
$c.Value = '0.0₃0769'
$c.ClearFormats()
$ws.Range('E16').Value = '  -0.14%  '
$c = $ws.Range('D17')
$c.NumberFormat = '@'
$c.Value = '62.76'
$c.ClearFormats()
$ws.Range('E17').Value = '  -0.81%  '
$c = $ws.Range('D18')
$c.NumberFormat = '@'
$c.Value = '25.801.47'
$c.ClearFormats()
$ws.Range('E18').Value = '  -0.29%  '
$ws.Range('E19').Value = '  -0.26%  '
$c = $ws.Range('D20')
$c.NumberFormat = '@'
$c.Value = '4.46'
$c.ClearFormats()
$ws.Range('E20').Value = '  +1.46%  '
$c = $ws.Range('D21')
$c.NumberFormat = '@'
$c.Value = '194.22'
$c.ClearFormats()
$ws.Range('E21').Value = '  +0.06%  '
$c = $ws.Range('D22')
$c.NumberFormat = '@'
$c.Value = '9.96'
$c.ClearFormats()
$ws.Range('E22').Value = '  -0.02%  '
$c = $ws.Range('D23')
$c.NumberFormat = '@'
$c.Value = '6.29'
$c.ClearFormats()
$ws.Range('E23').Value = '  +1.88%  '
$ws.Range('E24').Value = '  -0.21%  '
$ws.Range('E25').Value = '  +2.66%  '
$c = $ws.Range('D26')
$c.NumberFormat = '@'
$c.Value = '142.80'
$c.ClearFormats()
$ws.Range('E26').Value = '  +2.97%  '
$ws.Range('E27').Value = '  -0.09%  '
$ws.Range('E28').Value = '  +0.69%  '
$ws.Range('E29').Value = '  +0.09%  '
$ws.Range('E30').Value = '  -0.25%  '
$c = $ws.Range('D31')
$c.NumberFormat = '@'
$c.Value = '0.0493'
$c.ClearFormats()
$ws.Range('E31').Value = '  -0.79%  '
$ws.Range('E32').Value = '  +1.39%  '
$ws.Range('E33').Value = '  -0.27%  '
$ws.Range('E34').Value = '  +0.89%  '
$c = $ws.Range('D35')
$c.NumberFormat = '@'
$c.Value = '2.39'
$c.ClearFormats()
$ws.Range('E35').Value = '  -0.10%  '
$c = $ws.Range('D36')
$c.NumberFormat = '@'
$c.Value = '0.906'
$c.ClearFormats()
$ws.Range('E36').Value = '  +0.16%  '
$c = $ws.Range('D37')
$c.NumberFormat = '@'
$c.Value = '1.133.15'
$c.ClearFormats()
$ws.Range('E37').Value = '  -0.36%  '
$c = $ws.Range('D38')
$c.NumberFormat = '@'
$c.Value = '2.53'
$c.ClearFormats()
$ws.Range('E38').Value = '  -1.66%  '
$c = $ws.Range('D39')
$c.NumberFormat = '@'
$c.Value = '0.546'
$c.ClearFormats()
$ws.Range('E39').Value = '  -1.79%  '
$c = $ws.Range('D40')
$c.NumberFormat = '@'
$c.Value = '0.0157'
$c.ClearFormats()
$ws.Range('E41').Value = '  -0.36%  '
$ws.Range('E42').Value = '  +2.26%  '
$c = $ws.Range('D43')
$c.NumberFormat = '@'
$c.Value = '100.57'
$c.ClearFormats()
$ws.Range('E43').Value = '  +0.95%  '
$ws.Range('E44').Value = '  +0.79%  '
$c = $ws.Range('D45')
$c.NumberFormat = '@'
$c.Value = '1.768.66'
$c.ClearFormats()
$ws.Range('E45').Value = '  -0.45%  '
$ws.Range('E46').Value = '  +1.15%  '
$c = $ws.Range('D47')
$c.NumberFormat = '@'
$c.Value = '55.23'
$c.ClearFormats()
$ws.Range('E47').Value = '  -0.59%  '
$ws.Range('E48').Value = '  -0.36%  '
$c = $ws.Range('D49')
$c.NumberFormat = '@'
$c.Value = '0.416'
$c.ClearFormats()
$ws.Range('E49').Value = '  -2.31%  '
$ws.Range('E50').Value = '  +0.27%  '
$ws.Range('E51').Value = '  -3.35%  '
